$d = $word.ActiveDocument

$replacements = @(
    @("698×3=", "792×3="),
    @("209×4=", "743×6="),
    @("114×4=", "368×5="),
    @("437×2=", "537×6="),
    @("369×6=", "110×8="),
    @("543×8=", "212×6="),
    @("204×6=", "377×8="),
    @("178×2=", "908×6="),
    @("910×2=", "851×2="),
    @("400×3=", "448×4="),
    @("349×3=", "188×3="),
    @("206×2=", "466×7="),
    @("454×6=", "266×7="),
    @("923×9=", "576×8="),
    @("353×2=", "446×8="),
    @("315×6=", "743×4="),
    @("450×4=", "860×4="),
    @("418×3=", "195×2="),
    @("254×5=", "126×8="),
    @("981×7=", "586×3="),
    @("109×8=", "441×3="),
    @("230×6=", "613×6="),
    @("996×9=", "832×4="),
    @("678×9=", "824×6="),
    @("690×3=", "515×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
